# Update the raw metric figures on the "Metrics" sheet (B2:B13). Everything
# on the "today" sheet that reads from these cells (B11:B22 via
# "=Metrics!Bn", plus the dependent E11:E22/F11:F22 formulas) recalculates
# automatically once these source values change.
$wb = $excel.ActiveWorkbook
$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value = 285737.69000000006
$wsMetrics.Range("B3").Value = 251827.84000000003
$wsMetrics.Range("B4").Value = 88233.87000000001
$wsMetrics.Range("B5").Value = 11631
$wsMetrics.Range("B6").Value = 5081983.4400000013
$wsMetrics.Range("B7").Value = 4293904.5200000005
$wsMetrics.Range("B8").Value = 1495193.7000000002
$wsMetrics.Range("B9").Value = 197838
$wsMetrics.Range("B10").Value = 33547364.430000011
$wsMetrics.Range("B11").Value = 31569179.680000003
$wsMetrics.Range("B12").Value = 11776915.74
$wsMetrics.Range("B13").Value = 1295468

# Move the saved cursor/selection on "Metrics" to E15 (was E22).
[void]$wsMetrics.Range("E15").Select()

# Move the saved cursor/selection on "today" to F7 (was F8). Select it last
# so "today" stays the active/selected sheet tab, matching the workbook's
# state before the edit.
$wsToday = $wb.Worksheets.Item("today")
[void]$wsToday.Range("F7").Select()
